$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtering save games)
$data = @{
    2 = @(1.505614041169197, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0, 2.210719231951476)
    3 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 246.9852506941017, 0, 250.7519314561143)
    4 = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 0, 20.64246832346449)
    5 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    6 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
